$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 22

# Column A holds a date-like string ("12/16/2025") that must stay stored as
# literal text (matching the other rows' inlineStr date cells), not get
# auto-converted into a date serial number. Temporarily force a text number
# format so entry isn't reinterpreted, then restore the default "Normal"
# style afterwards so the cell ends up with no explicit style (same as the
# rest of the data rows).
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "12/16/2025"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = 11757.73
$ws.Cells.Item($row, 3).Value = 0.2092236483795707
$ws.Cells.Item($row, 4).Value = 0.7907763516204293
$ws.Cells.Item($row, 5).Value = -148.26
$ws.Cells.Item($row, 6).Value = -32.23
$ws.Cells.Item($row, 7).Value = -21219.02
$ws.Cells.Item($row, 8).Value = -69.53
$ws.Cells.Item($row, 9).Value = -492.85
$ws.Cells.Item($row, 10).Value = -16.69
